$wb = $excel.ActiveWorkbook

# Hospital: closing-schedule shifted (2nd marking algorithm)
$ws = $wb.Worksheets.Item("Hospital")
$ws.Range("B9").Value  = "livre"
$ws.Range("B14").Value = "fechado"
$ws.Range("B21").Value = "fechado"
$ws.Range("B33").Value = "livre"
$ws.Range("B34").Value = "livre"
$ws.Range("B38").Value = "fechado"
$ws.Range("B45").Value = "fechado"
$ws.Range("B46").Value = "fechado"
$ws.Range("B57").Value = "livre"
$ws.Range("B58").Value = "livre"
$ws.Range("B59").Value = "livre"
$ws.Range("B62").Value = "fechado"
$ws.Range("B69").Value = "fechado"
$ws.Range("B70").Value = "fechado"
$ws.Range("B71").Value = "fechado"
$ws.Range("B81").Value = "livre"
$ws.Range("B82").Value = "livre"
$ws.Range("B83").Value = "livre"
$ws.Range("B86").Value = "fechado"
$ws.Range("B93").Value = "fechado"
$ws.Range("B94").Value = "fechado"
$ws.Range("B95").Value = "fechado"
$ws.Range("B96").Value = "fechado"

# Patient1: "marcado"/"ocupador" statuses consolidated into "ocupado"
$ws = $wb.Worksheets.Item("Patient1")
$ws.Range("B25").Value = "ocupado"
$ws.Range("B44").Value = "ocupado"
$ws.Range("B45").Value = "ocupado"
$ws.Range("B46").Value = "ocupado"
$ws.Range("B47").Value = "ocupado"
$ws.Range("B48").Value = "ocupado"
$ws.Range("B49").Value = "ocupado"
$ws.Range("B59").Value = "ocupado"

# Patient3: "marcado" status consolidated into "ocupado"
$ws = $wb.Worksheets.Item("Patient3")
$ws.Range("B23").Value = "ocupado"
